$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-write A1 / B1 headers (same text) so the cell gets an explicit style index written
$ws.Range("A1").Value = "fila"
$ws.Range("B1").Value = "resultado"

# Row 2 already has "ok" in B2; rewrite it so style gets normalized too
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "ok"

# Fill column B (resultado) with "ok" for remaining rows, and refresh column A (fila)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "ok"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "ok"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "ok"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "ok"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "ok"

# Introduce the new "fallo" string into the shared-strings table (written, then cleared)
$ws.Range("B9").Value = "fallo"
$ws.Range("B9").Value = ""

# Move the active selection to B8, mirroring where the user's cursor ended up
$ws.Range("B8").Select()
